# Auto-generated Excel COM-interop script
# Applies numeric updates to the FFXIV leve profit tables across all 8 class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 700
$ws.Range("I6").Value = 700
$ws.Range("K6").Value = 2100
$ws.Range("M6").Value = -1988
$ws.Range("H9").Value = 672.8182
$ws.Range("I9").Value = 705
$ws.Range("J9").Value = 351
$ws.Range("K9").Value = 705
$ws.Range("L9").Value = 351
$ws.Range("M9").Value = -536
$ws.Range("N9").Value = -689
$ws.Range("H12").Value = 138.85715
$ws.Range("I12").Value = 138.85715
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 138.85715
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 31.14285000000001
$ws.Range("H21").Value = 8312.5
$ws.Range("J21").Value = 8312.5
$ws.Range("L21").Value = 8312.5
$ws.Range("N21").Value = -9248.5
$ws.Range("H23").Value = 8312.5
$ws.Range("J23").Value = 8312.5
$ws.Range("L23").Value = 8312.5
$ws.Range("N23").Value = -8780.5
$ws.Range("H29").Value = 1506.5834
$ws.Range("J29").Value = 2485.7144
$ws.Range("L29").Value = 7457.1432
$ws.Range("N29").Value = -8019.1432
$ws.Range("H38").Value = 1242.25
$ws.Range("I38").Value = 50.875
$ws.Range("J38").Value = 3625
$ws.Range("K38").Value = 152.625
$ws.Range("L38").Value = 10875
$ws.Range("M38").Value = 219.375
$ws.Range("N38").Value = -11619
$ws.Range("H58").Value = 1006.8182
$ws.Range("I58").Value = 586.1111
$ws.Range("J58").Value = 2900
$ws.Range("K58").Value = 1758.3333
$ws.Range("L58").Value = 8700
$ws.Range("M58").Value = -1608.3333
$ws.Range("N58").Value = -9000
$ws.Range("H61").Value = 174.54546
$ws.Range("I61").Value = 94
$ws.Range("K61").Value = 282
$ws.Range("M61").Value = -110
$ws.Range("H64").Value = 3038.4814
$ws.Range("I64").Value = 2979.0908
$ws.Range("J64").Value = 3079.3125
$ws.Range("K64").Value = 2979.0908
$ws.Range("L64").Value = 3079.3125
$ws.Range("M64").Value = -2731.0908
$ws.Range("N64").Value = -3575.3125
$ws.Range("H67").Value = 3038.4814
$ws.Range("I67").Value = 2979.0908
$ws.Range("J67").Value = 3079.3125
$ws.Range("K67").Value = 2979.0908
$ws.Range("L67").Value = 3079.3125
$ws.Range("M67").Value = -2121.0908
$ws.Range("N67").Value = -4795.3125
$ws.Range("H76").Value = 3309.2856
$ws.Range("I76").Value = 2999.9375
$ws.Range("K76").Value = 2999.9375
$ws.Range("M76").Value = -2684.9375
$ws.Range("H79").Value = 3309.2856
$ws.Range("I79").Value = 2999.9375
$ws.Range("K79").Value = 2999.9375
$ws.Range("M79").Value = -1907.9375
$ws.Range("H87").Value = 25790
$ws.Range("J87").Value = 25790
$ws.Range("L87").Value = 25790
$ws.Range("N87").Value = -28286
$ws.Range("H90").Value = 25790
$ws.Range("J90").Value = 25790
$ws.Range("L90").Value = 77370
$ws.Range("N90").Value = -89850
$ws.Range("H138").Value = 3386.5942
$ws.Range("I138").Value = 1853.6111
$ws.Range("J138").Value = 3927.647
$ws.Range("K138").Value = 5560.8333
$ws.Range("L138").Value = 11782.941
$ws.Range("M138").Value = -420.8333000000002
$ws.Range("N138").Value = -22062.941
$ws.Range("N12").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5229.5386
$ws.Range("I31").Value = 5229.5386
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5229.5386
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4935.5386
$ws.Range("H32").Value = 8163.29
$ws.Range("I32").Value = 5952.4805
$ws.Range("J32").Value = 15564.695
$ws.Range("K32").Value = 5952.4805
$ws.Range("L32").Value = 15564.695
$ws.Range("M32").Value = -5665.4805
$ws.Range("N32").Value = -16138.695
$ws.Range("H63").Value = 3922.1428
$ws.Range("I63").Value = 2691
$ws.Range("K63").Value = 2691
$ws.Range("M63").Value = -2005
$ws.Range("H66").Value = 3922.1428
$ws.Range("I66").Value = 2691
$ws.Range("K66").Value = 13455
$ws.Range("M66").Value = -10023
$ws.Range("H88").Value = 166801680
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 200161420
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 200161420
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -200162232
$ws.Range("H91").Value = 166801680
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 200161420
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 200161420
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -200164228
$ws.Range("N31").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 18657
$ws.Range("J103").Value = 18657
$ws.Range("L103").Value = 18657
$ws.Range("N103").Value = -21001
$ws.Range("H132").Value = 43450
$ws.Range("J132").Value = 43450
$ws.Range("L132").Value = 43450
$ws.Range("N132").Value = -53570
$ws.Range("H134").Value = 313835.2
$ws.Range("I134").Value = 371243.03
$ws.Range("J134").Value = 3832.8
$ws.Range("K134").Value = 1113729.09
$ws.Range("L134").Value = 11498.4
$ws.Range("M134").Value = -1111194.09
$ws.Range("N134").Value = -16568.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3438.4443
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 3580.75
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 3580.75
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -4828.75
$ws.Range("H65").Value = 3438.4443
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 3580.75
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 17903.75
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -24143.75
$ws.Range("H99").Value = 4158.2666
$ws.Range("I99").Value = 4337.4
$ws.Range("J99").Value = 3800
$ws.Range("K99").Value = 4337.4
$ws.Range("L99").Value = 3800
$ws.Range("M99").Value = -2839.4
$ws.Range("N99").Value = -6796
$ws.Range("H126").Value = 4158.2666
$ws.Range("I126").Value = 4337.4
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 13012.2
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -10542.2
$ws.Range("N126").Value = -16340
$ws.Range("H141").Value = 43706
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 43706
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 43706
$ws.Range("N141").Value = -54066
$ws.Range("M141").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1120.6364
$ws.Range("I7").Value = 189.8
$ws.Range("J7").Value = 1896.3334
$ws.Range("K7").Value = 569.4000000000001
$ws.Range("L7").Value = 5689.0002
$ws.Range("M7").Value = -457.4000000000001
$ws.Range("N7").Value = -5913.0002
$ws.Range("H131").Value = 1551198.9
$ws.Range("I131").Value = 5555945.5
$ws.Range("J131").Value = 974.4516
$ws.Range("K131").Value = 16667836.5
$ws.Range("L131").Value = 2923.3548
$ws.Range("M131").Value = -16662796.5
$ws.Range("N131").Value = -13003.3548
$ws.Range("H139").Value = 7874.8237
$ws.Range("I139").Value = 1776.5714
$ws.Range("J139").Value = 36333.332
$ws.Range("K139").Value = 5329.7142
$ws.Range("L139").Value = 108999.996
$ws.Range("M139").Value = -189.7142000000003
$ws.Range("N139").Value = -119279.996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 69122280
$ws.Range("I70").Value = 207354640
$ws.Range("J70").Value = 6100
$ws.Range("K70").Value = 207354640
$ws.Range("L70").Value = 6100
$ws.Range("M70").Value = -207354370
$ws.Range("N70").Value = -6640
$ws.Range("H73").Value = 69122280
$ws.Range("I73").Value = 207354640
$ws.Range("J73").Value = 6100
$ws.Range("K73").Value = 207354640
$ws.Range("L73").Value = 6100
$ws.Range("M73").Value = -207353704
$ws.Range("N73").Value = -7972
$ws.Range("H80").Value = 3652.1428
$ws.Range("I80").Value = 5667.5
$ws.Range("J80").Value = 2846
$ws.Range("K80").Value = 5667.5
$ws.Range("L80").Value = 2846
$ws.Range("M80").Value = -4669.5
$ws.Range("N80").Value = -4842
$ws.Range("H83").Value = 3652.1428
$ws.Range("I83").Value = 5667.5
$ws.Range("J83").Value = 2846
$ws.Range("K83").Value = 28337.5
$ws.Range("L83").Value = 14230
$ws.Range("M83").Value = -23345.5
$ws.Range("N83").Value = -24214
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
$ws.Range("H132").Value = 3092
$ws.Range("I132").Value = 2279
$ws.Range("J132").Value = 4243.75
$ws.Range("K132").Value = 6837
$ws.Range("L132").Value = 12731.25
$ws.Range("M132").Value = -4307
$ws.Range("N132").Value = -17791.25
$ws.Range("H134").Value = 33725.2
$ws.Range("J134").Value = 33725.2
$ws.Range("L134").Value = 101175.6
$ws.Range("N134").Value = -106245.6
$ws.Range("H136").Value = 26484.902
$ws.Range("J136").Value = 26484.902
$ws.Range("L136").Value = 79454.70599999999
$ws.Range("N136").Value = -84554.70599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H101").Value = 34700
$ws.Range("J101").Value = 34700
$ws.Range("L101").Value = 34700
$ws.Range("N101").Value = -41190
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1391.9166
$ws.Range("I122").Value = 988.94446
$ws.Range("K122").Value = 2966.83338
$ws.Range("M122").Value = -516.83338
